# "Sign in and register feature completed"
#
# Updates the DSAlgo BDD framework test-data workbook:
#   - loginData: refresh username/password test values, add a hyperlink
#     on the password cell, move the saved selection to A4
#   - pythonCode: tidy up the expected "result" header/value, become the
#     active sheet/tab with selection B15
#   - registerData: no longer the active tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# loginData sheet
# ---------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("loginData")

$wsLogin.Range("A2").Value = "sdetwarriors"
$wsLogin.Range("B3").Value = "sdet@146"
$wsLogin.Range("A4").Value = "n2324435"

# New hyperlink on the password cell (B3), displaying "sdet@146"
# (Excel's AutoFormat turns "@"-containing text typed into a cell into a
# mailto: link; that's how this hyperlink came to exist on a password cell)
$wsLogin.Hyperlinks.Add($wsLogin.Range("B3"), "mailto:sdet@146", "", "", "sdet@146")

# Restore the saved selection to A4 (without making this sheet active)
$wsLogin.Range("A4").Select()

# ---------------------------------------------------------------------
# pythonCode sheet
# ---------------------------------------------------------------------
$wsPython = $wb.Worksheets.Item("pythonCode")

$wsPython.Range("B1").Value = "result"
$wsPython.Range("A2").Value = 'print("hello")'

# ---------------------------------------------------------------------
# registerData sheet (was the active tab; stays at its own selection)
# ---------------------------------------------------------------------
$wsRegister = $wb.Worksheets.Item("registerData")
$wsRegister.Range("C5").Select()

# ---------------------------------------------------------------------
# pythonCode becomes the active sheet/tab, selection -> B15
# ---------------------------------------------------------------------
$wsPython.Activate()
$wsPython.Range("B15").Select()
